$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.241344451904297
$ws.Range("B1").Value = 2.356509923934937
$ws.Range("C1").Value = 3.712232351303101
$ws.Range("D1").Value = 3.422984838485718
$ws.Range("E1").Value = 1.203454256057739
